# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Leve profit-tracking sheets
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns H-N)
# as refreshed by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# ALC row 2 (Leve Item ID 5489)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2,8).Value = 5377.857  # H2: 4805.625 -> 5377.857
$ws.Cells.Item(2,9).Value = 2850  # I2: 2166.6667 -> 2850
$ws.Cells.Item(2,11).Value = 2850  # K2: 2166.6667 -> 2850
$ws.Cells.Item(2,13).Value = -2737  # M2: -2053.6667 -> -2737

# ALC row 6 (Leve Item ID 4564)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6,8).Value = 789.8823  # H6: 746.05554 -> 789.8823
$ws.Cells.Item(6,9).Value = 321.5  # I6: 275.7143 -> 321.5
$ws.Cells.Item(6,11).Value = 964.5  # K6: 827.1428999999999 -> 964.5
$ws.Cells.Item(6,13).Value = -852.5  # M6: -715.1428999999999 -> -852.5

# ALC row 28 (Leve Item ID 27772)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28,8).Value = 264.83334  # H28: 263.16666 -> 264.83334
$ws.Cells.Item(28,9).Value = 190  # I28: 180 -> 190
$ws.Cells.Item(28,11).Value = 190  # K28: 180 -> 190
$ws.Cells.Item(28,13).Value = 295  # M28: 305 -> 295

# ALC row 70 (Leve Item ID 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70,8).Value = 2188.8  # H70: 1496.1 -> 2188.8
$ws.Cells.Item(70,9).Value = 0  # I70: 810 -> 0
$ws.Cells.Item(70,10).Value = 2188.8  # J70: 1953.5 -> 2188.8
$ws.Cells.Item(70,11).Value = 0  # K70: 2430 -> 0
$ws.Cells.Item(70,12).Value = 6566.400000000001  # L70: 5860.5 -> 6566.400000000001
$ws.Cells.Item(70,13).ClearContents()  # M70: -2160 -> (removed)
$ws.Cells.Item(70,14).Value = -7106.400000000001  # N70: -6400.5 -> -7106.400000000001

# ALC row 73 (Leve Item ID 12604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73,8).Value = 2188.8  # H73: 1496.1 -> 2188.8
$ws.Cells.Item(73,9).Value = 0  # I73: 810 -> 0
$ws.Cells.Item(73,10).Value = 2188.8  # J73: 1953.5 -> 2188.8
$ws.Cells.Item(73,11).Value = 0  # K73: 2430 -> 0
$ws.Cells.Item(73,12).Value = 6566.400000000001  # L73: 5860.5 -> 6566.400000000001
$ws.Cells.Item(73,13).ClearContents()  # M73: -1494 -> (removed)
$ws.Cells.Item(73,14).Value = -8438.400000000001  # N73: -7732.5 -> -8438.400000000001

# ALC row 113 (Leve Item ID 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113,8).Value = 3726.4285  # H113: 3743.8462 -> 3726.4285
$ws.Cells.Item(113,10).Value = 4167  # J113: 4300.4 -> 4167
$ws.Cells.Item(113,12).Value = 4167  # L113: 4300.4 -> 4167
$ws.Cells.Item(113,14).Value = -10675  # N113: -10808.4 -> -10675

# ALC row 115 (Leve Item ID 27957)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(115,8).Value = 1317.3334  # H115: 1263.6 -> 1317.3334
$ws.Cells.Item(115,9).Value = 869.4  # I115: 854.3333 -> 869.4
$ws.Cells.Item(115,10).Value = 1877.25  # J115: 1877.5 -> 1877.25
$ws.Cells.Item(115,11).Value = 2608.2  # K115: 2562.9999 -> 2608.2
$ws.Cells.Item(115,12).Value = 5631.75  # L115: 5632.5 -> 5631.75
$ws.Cells.Item(115,13).Value = -1041.2  # M115: -995.9998999999998 -> -1041.2
$ws.Cells.Item(115,14).Value = -8765.75  # N115: -8766.5 -> -8765.75

# ALC row 126 (Leve Item ID 34391)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(126,8).Value = 0  # H126: 139995 -> 0
$ws.Cells.Item(126,10).Value = 0  # J126: 139995 -> 0
$ws.Cells.Item(126,12).Value = 0  # L126: 139995 -> 0
$ws.Cells.Item(126,14).ClearContents()  # N126: -149875 -> (removed)

# ALC row 128 (Leve Item ID 34540)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(128,8).Value = 124434.11  # H128: 156656.67 -> 124434.11
$ws.Cells.Item(128,10).Value = 124434.11  # J128: 156656.67 -> 124434.11
$ws.Cells.Item(128,12).Value = 124434.11  # L128: 156656.67 -> 124434.11
$ws.Cells.Item(128,14).Value = -134394.11  # N128: -166616.67 -> -134394.11

# ALC row 129 (Leve Item ID 36115)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(129,8).Value = 1744.125  # H129: 1618.5 -> 1744.125
$ws.Cells.Item(129,9).Value = 1075.5  # I129: 1135.4286 -> 1075.5
$ws.Cells.Item(129,10).Value = 3750  # J129: 5000 -> 3750
$ws.Cells.Item(129,11).Value = 3226.5  # K129: 3406.2858 -> 3226.5
$ws.Cells.Item(129,12).Value = 11250  # L129: 15000 -> 11250
$ws.Cells.Item(129,13).Value = 1773.5  # M129: 1593.7142 -> 1773.5
$ws.Cells.Item(129,14).Value = -21250  # N129: -25000 -> -21250

# ALC row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132,8).Value = 1520.3684  # H132: 1544.5676 -> 1520.3684
$ws.Cells.Item(132,9).Value = 1312.4445  # I132: 1338.8846 -> 1312.4445
$ws.Cells.Item(132,11).Value = 3937.3335  # K132: 4016.6538 -> 3937.3335
$ws.Cells.Item(132,13).Value = -1407.3335  # M132: -1486.6538 -> -1407.3335

# ARM row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132,8).Value = 1898.1163  # H132: 1794.9302 -> 1898.1163
$ws.Cells.Item(132,9).Value = 1578.0834  # I132: 1549.2703 -> 1578.0834
$ws.Cells.Item(132,10).Value = 3544  # J132: 3309.8333 -> 3544
$ws.Cells.Item(132,11).Value = 4734.2502  # K132: 4647.810899999999 -> 4734.2502
$ws.Cells.Item(132,12).Value = 10632  # L132: 9929.499899999999 -> 10632
$ws.Cells.Item(132,13).Value = -2204.2502  # M132: -2117.810899999999 -> -2204.2502
$ws.Cells.Item(132,14).Value = -15692  # N132: -14989.4999 -> -15692

# BSM row 107 (Leve Item ID 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107,8).Value = 14581.765  # H107: 14668.823 -> 14581.765
$ws.Cells.Item(107,9).Value = 18118.46  # I107: 18232.309 -> 18118.46
$ws.Cells.Item(107,11).Value = 18118.46  # K107: 18232.309 -> 18118.46
$ws.Cells.Item(107,13).Value = -16198.46  # M107: -16312.309 -> -16198.46

# CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 2746.8708  # H31: 2118.0408 -> 2746.8708
$ws.Cells.Item(31,9).Value = 5833.3335  # I31: 1683.2632 -> 5833.3335
$ws.Cells.Item(31,10).Value = 2416.1785  # J31: 2393.4 -> 2416.1785
$ws.Cells.Item(31,11).Value = 5833.3335  # K31: 1683.2632 -> 5833.3335
$ws.Cells.Item(31,12).Value = 2416.1785  # L31: 2393.4 -> 2416.1785
$ws.Cells.Item(31,13).Value = -5538.3335  # M31: -1388.2632 -> -5538.3335
$ws.Cells.Item(31,14).Value = -3006.1785  # N31: -2983.4 -> -3006.1785

# CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34,8).Value = 2746.8708  # H34: 2118.0408 -> 2746.8708
$ws.Cells.Item(34,9).Value = 5833.3335  # I34: 1683.2632 -> 5833.3335
$ws.Cells.Item(34,10).Value = 2416.1785  # J34: 2393.4 -> 2416.1785
$ws.Cells.Item(34,11).Value = 5833.3335  # K34: 1683.2632 -> 5833.3335
$ws.Cells.Item(34,12).Value = 2416.1785  # L34: 2393.4 -> 2416.1785
$ws.Cells.Item(34,13).Value = -5631.3335  # M34: -1481.2632 -> -5631.3335
$ws.Cells.Item(34,14).Value = -2820.1785  # N34: -2797.4 -> -2820.1785

# CRP row 122 (Leve Item ID 36196)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122,8).Value = 1575.875  # H122: 1627.9333 -> 1575.875
$ws.Cells.Item(122,9).Value = 1450.0834  # I122: 1509.6364 -> 1450.0834
$ws.Cells.Item(122,11).Value = 4350.2502  # K122: 4528.9092 -> 4350.2502
$ws.Cells.Item(122,13).Value = -1900.2502  # M122: -2078.9092 -> -1900.2502

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134,8).Value = 1667.5588  # H134: 1704.4546 -> 1667.5588
$ws.Cells.Item(134,9).Value = 1550.5416  # I134: 1598.3914 -> 1550.5416
$ws.Cells.Item(134,11).Value = 4651.6248  # K134: 4795.174199999999 -> 4651.6248
$ws.Cells.Item(134,13).Value = -2116.6248  # M134: -2260.174199999999 -> -2116.6248

# CUL row 38 (Leve Item ID 4860)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38,8).Value = 409.57144  # H38: 468.66666 -> 409.57144
$ws.Cells.Item(38,9).Value = 567.5  # I38: 738.3333 -> 567.5
$ws.Cells.Item(38,11).Value = 1702.5  # K38: 2214.9999 -> 1702.5
$ws.Cells.Item(38,13).Value = -1355.5  # M38: -1867.9999 -> -1355.5

# CUL row 46 (Leve Item ID 4701)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46,8).Value = 33333372  # H46: 16666900 -> 33333372
$ws.Cells.Item(46,10).Value = 0  # J46: 428.33334 -> 0
$ws.Cells.Item(46,12).Value = 0  # L46: 1285.00002 -> 0
$ws.Cells.Item(46,14).ClearContents()  # N46: -1467.00002 -> (removed)

# CUL row 58 (Leve Item ID 4703)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(58,8).Value = 14166.111  # H58: 10499.667 -> 14166.111
$ws.Cells.Item(58,10).Value = 17142.857  # J58: 11538.462 -> 17142.857
$ws.Cells.Item(58,12).Value = 51428.571  # L58: 34615.386 -> 51428.571
$ws.Cells.Item(58,14).Value = -51684.571  # N58: -34871.386 -> -51684.571

# CUL row 98 (Leve Item ID 19843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98,8).Value = 1501.1818  # H98: 1510.9 -> 1501.1818
$ws.Cells.Item(98,10).Value = 1447  # J98: 1461.3334 -> 1447
$ws.Cells.Item(98,12).Value = 4341  # L98: 4384.0002 -> 4341
$ws.Cells.Item(98,14).Value = -7337  # N98: -7380.0002 -> -7337

# CUL row 137 (Leve Item ID 44088)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137,8).Value = 1474.8572  # H137: 1474.1428 -> 1474.8572
$ws.Cells.Item(137,10).Value = 2595  # J137: 2593.125 -> 2595
$ws.Cells.Item(137,12).Value = 7785  # L137: 7779.375 -> 7785
$ws.Cells.Item(137,14).Value = -17985  # N137: -17979.375 -> -17985

# GSM row 70 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70,8).Value = 6489.3335  # H70: 5766 -> 6489.3335
$ws.Cells.Item(70,9).Value = 6284  # I70: 5740.75 -> 6284
$ws.Cells.Item(70,10).Value = 6900  # J70: 5799.6665 -> 6900
$ws.Cells.Item(70,11).Value = 6284  # K70: 5740.75 -> 6284
$ws.Cells.Item(70,12).Value = 6900  # L70: 5799.6665 -> 6900
$ws.Cells.Item(70,13).Value = -6014  # M70: -5470.75 -> -6014
$ws.Cells.Item(70,14).Value = -7440  # N70: -6339.6665 -> -7440

# GSM row 73 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73,8).Value = 6489.3335  # H73: 5766 -> 6489.3335
$ws.Cells.Item(73,9).Value = 6284  # I73: 5740.75 -> 6284
$ws.Cells.Item(73,10).Value = 6900  # J73: 5799.6665 -> 6900
$ws.Cells.Item(73,11).Value = 6284  # K73: 5740.75 -> 6284
$ws.Cells.Item(73,12).Value = 6900  # L73: 5799.6665 -> 6900
$ws.Cells.Item(73,13).Value = -5348  # M73: -4804.75 -> -5348
$ws.Cells.Item(73,14).Value = -8772  # N73: -7671.6665 -> -8772

# GSM row 92 (Leve Item ID 18094)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(92,8).Value = 9127.75  # H92: 12537 -> 9127.75
$ws.Cells.Item(92,10).Value = 9127.75  # J92: 12537 -> 9127.75
$ws.Cells.Item(92,12).Value = 9127.75  # L92: 12537 -> 9127.75
$ws.Cells.Item(92,14).Value = -12871.75  # N92: -16281 -> -12871.75

# GSM row 107 (Leve Item ID 27802)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107,8).Value = 49334.81  # H107: 47146.816 -> 49334.81
$ws.Cells.Item(107,9).Value = 126093.5  # I107: 112216.336 -> 126093.5
$ws.Cells.Item(107,11).Value = 126093.5  # K107: 112216.336 -> 126093.5
$ws.Cells.Item(107,13).Value = -124173.5  # M107: -110296.336 -> -124173.5

# GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132,8).Value = 5517420.5  # H132: 5409249.5 -> 5517420.5
$ws.Cells.Item(132,9).Value = 5219.5127  # I132: 5106.525 -> 5219.5127
$ws.Cells.Item(132,11).Value = 15658.5381  # K132: 15319.575 -> 15658.5381
$ws.Cells.Item(132,13).Value = -13128.5381  # M132: -12789.575 -> -13128.5381

# LTW row 7 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value = 8458.762000000001  # H7: 8606.237999999999 -> 8458.762000000001
$ws.Cells.Item(7,9).Value = 4703.5  # I7: 4961.5835 -> 4703.5
$ws.Cells.Item(7,11).Value = 4703.5  # K7: 4961.5835 -> 4703.5
$ws.Cells.Item(7,13).Value = -4591.5  # M7: -4849.5835 -> -4591.5

# LTW row 22 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22,8).Value = 7281.4346  # H22: 7921.4707 -> 7281.4346
$ws.Cells.Item(22,9).Value = 6274.2  # I22: 6477.5557 -> 6274.2
$ws.Cells.Item(22,10).Value = 8056.231  # J22: 9545.875 -> 8056.231
$ws.Cells.Item(22,11).Value = 6274.2  # K22: 6477.5557 -> 6274.2
$ws.Cells.Item(22,12).Value = 8056.231  # L22: 9545.875 -> 8056.231
$ws.Cells.Item(22,13).Value = -5979.2  # M22: -6182.5557 -> -5979.2
$ws.Cells.Item(22,14).Value = -8646.231  # N22: -10135.875 -> -8646.231

# LTW row 27 (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27,8).Value = 7281.4346  # H27: 7921.4707 -> 7281.4346
$ws.Cells.Item(27,9).Value = 6274.2  # I27: 6477.5557 -> 6274.2
$ws.Cells.Item(27,10).Value = 8056.231  # J27: 9545.875 -> 8056.231
$ws.Cells.Item(27,11).Value = 6274.2  # K27: 6477.5557 -> 6274.2
$ws.Cells.Item(27,12).Value = 8056.231  # L27: 9545.875 -> 8056.231
$ws.Cells.Item(27,13).Value = -6167.2  # M27: -6370.5557 -> -6167.2
$ws.Cells.Item(27,14).Value = -8270.231  # N27: -9759.875 -> -8270.231

# LTW row 46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46,8).Value = 7940.5  # H46: 7033.421 -> 7940.5
$ws.Cells.Item(46,9).Value = 8653.786  # I46: 7934.1875 -> 8653.786
$ws.Cells.Item(46,10).Value = 2947.5  # J46: 2229.3333 -> 2947.5
$ws.Cells.Item(46,11).Value = 8653.786  # K46: 7934.1875 -> 8653.786
$ws.Cells.Item(46,12).Value = 2947.5  # L46: 2229.3333 -> 2947.5
$ws.Cells.Item(46,13).Value = -8465.786  # M46: -7746.1875 -> -8465.786
$ws.Cells.Item(46,14).Value = -3323.5  # N46: -2605.3333 -> -3323.5

# LTW row 47 (Leve Item ID 3138)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(47,8).Value = 29999  # H47: 0 -> 29999
$ws.Cells.Item(47,10).Value = 29999  # J47: 0 -> 29999
$ws.Cells.Item(47,12).Value = 29999  # L47: 0 -> 29999
$ws.Cells.Item(47,14).Value = -30979  # N47: None -> -30979

# LTW row 50 (Leve Item ID 3426)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(50,8).Value = 30000  # H50: 0 -> 30000
$ws.Cells.Item(50,9).Value = 30000  # I50: 0 -> 30000
$ws.Cells.Item(50,11).Value = 30000  # K50: 0 -> 30000
$ws.Cells.Item(50,13).Value = -29363  # M50: None -> -29363

# LTW row 51 (Leve Item ID 3423)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(51,8).Value = 26000  # H51: 0 -> 26000
$ws.Cells.Item(51,10).Value = 26000  # J51: 0 -> 26000
$ws.Cells.Item(51,12).Value = 26000  # L51: 0 -> 26000
$ws.Cells.Item(51,14).Value = -26956  # N51: None -> -26956

# LTW row 52 (Leve Item ID 3138)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(52,8).Value = 29999  # H52: 0 -> 29999
$ws.Cells.Item(52,10).Value = 29999  # J52: 0 -> 29999
$ws.Cells.Item(52,12).Value = 29999  # L52: 0 -> 29999
$ws.Cells.Item(52,14).Value = -30465  # N52: None -> -30465

# LTW row 124 (Leve Item ID 34264)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(124,8).Value = 49444.25  # H124: 77777 -> 49444.25
$ws.Cells.Item(124,10).Value = 49444.25  # J124: 77777 -> 49444.25
$ws.Cells.Item(124,12).Value = 49444.25  # L124: 77777 -> 49444.25
$ws.Cells.Item(124,14).Value = -59264.25  # N124: -87597 -> -59264.25

# LTW row 126 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126,8).Value = 8458.762000000001  # H126: 8606.237999999999 -> 8458.762000000001
$ws.Cells.Item(126,9).Value = 4703.5  # I126: 4961.5835 -> 4703.5
$ws.Cells.Item(126,11).Value = 14110.5  # K126: 14884.7505 -> 14110.5
$ws.Cells.Item(126,13).Value = -11640.5  # M126: -12414.7505 -> -11640.5

# LTW row 127 (Leve Item ID 34401)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(127,8).Value = 0  # H127: 81950 -> 0
$ws.Cells.Item(127,10).Value = 0  # J127: 81950 -> 0
$ws.Cells.Item(127,12).Value = 0  # L127: 81950 -> 0
$ws.Cells.Item(127,14).ClearContents()  # N127: -91870 -> (removed)

# LTW row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136,8).Value = 2090.158  # H136: 2062.1794 -> 2090.158
$ws.Cells.Item(136,10).Value = 2374.2  # J136: 2288.25 -> 2374.2
$ws.Cells.Item(136,12).Value = 7122.599999999999  # L136: 6864.75 -> 7122.599999999999
$ws.Cells.Item(136,14).Value = -12222.6  # N136: -11964.75 -> -12222.6

# WVR row 119 (Leve Item ID 26289)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119,8).Value = 0  # H119: 41000 -> 0
$ws.Cells.Item(119,10).Value = 0  # J119: 41000 -> 0
$ws.Cells.Item(119,12).Value = 0  # L119: 41000 -> 0
$ws.Cells.Item(119,14).ClearContents()  # N119: -50676 -> (removed)

# WVR row 122 (Leve Item ID 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122,8).Value = 3107.2222  # H122: 3205.0588 -> 3107.2222
$ws.Cells.Item(122,9).Value = 2580.077  # I122: 2674.75 -> 2580.077
$ws.Cells.Item(122,11).Value = 7740.231000000001  # K122: 8024.25 -> 7740.231000000001
$ws.Cells.Item(122,13).Value = -5290.231000000001  # M122: -5574.25 -> -5290.231000000001

# WVR row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132,8).Value = 1101.4706  # H132: 1155.8438 -> 1101.4706
$ws.Cells.Item(132,9).Value = 752.11536  # I132: 795.5 -> 752.11536
$ws.Cells.Item(132,11).Value = 2256.34608  # K132: 2386.5 -> 2256.34608
$ws.Cells.Item(132,13).Value = 273.6539199999997  # M132: 143.5 -> 273.6539199999997
